# Edit: merge the "Team #7" / " CSS Report" title runs into a single run,
# and move the "_GoBack" bookmark from the (empty) first paragraph down to
# the end of the title paragraph (right after the merged run's text, before
# its paragraph mark).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the existing "_GoBack" bookmark (currently sitting alone
# in the first, empty paragraph).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# ---------------------------------------------------------------------
# Step 2: merge the two title runs ("Team #7" + " CSS Report") into one
# run containing "Team #7 CSS Report". Both runs already share identical
# formatting, so replacing the whole phrase with itself collapses them
# into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Team #7 CSS Report", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Team #7 CSS Report", 2)

# ---------------------------------------------------------------------
# Step 3: re-insert the "_GoBack" bookmark at the end of the title
# paragraph (after the merged run's text, before the paragraph mark).
#
# Note: adding a bookmark directly to a zero-length Range placed exactly
# before a paragraph mark is unreliable here, so a one-character "anchor"
# is inserted, the bookmark is attached to the range spanning that
# character, and the anchor character is removed afterwards - leaving the
# bookmark collapsed at the correct position.
# ---------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Team #7 CSS Report*") {
        $titlePara = $p
        break
    }
}

$insertPos = $titlePara.Range.End - 1
$anchorRange = $d.Range($insertPos, $insertPos)
$anchorRange.InsertAfter("X")

$bookmarkRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$cleanupRange = $d.Range($insertPos, $insertPos + 1)
$cleanupRange.Text = ""
